# Applies the "Error Calculations and Plots" edit:
#  - Fill in / clear several missing-data cells in rows 2-25 (columns E/F)
#  - Remove two rows (RM 232 and SC 92) which shifts the remaining rows up,
#    and backfill previously-missing cells revealed in the shifted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell-level fixes within rows 2-25 (no row shift here) ---
$ws.Range("E3").Value = -5.7
$ws.Range("F4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F9").Value = 17.26
$ws.Range("F10").Value = 16.43
$ws.Range("F11").Value = 17.65
$ws.Range("F12").Value = 17.45
$ws.Range("F15").ClearContents()
$ws.Range("F17").ClearContents()
$ws.Range("F18").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("E23").ClearContents()

# --- Remove the "RM 232" row (row 26) entirely; rows below shift up by one ---
$ws.Rows("26").Delete()

# --- Remove the "SC 92" row (now row 27 after the previous shift) ---
$ws.Rows("27").Delete()

# --- Backfill values that were previously missing, now revealed after shift ---
# (former row 33 "SC 132" -> now row 31; former row 34 "SC 193" -> now row 32)
$ws.Range("F31").Value = 17.18
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39
